# The underlying commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml (bound to the slide master / "Design")  <->
#   ppt/theme/theme2.xml (bound to the notes master)
# Before:  theme1.xml = "Integral" colours, theme2.xml = "Office Theme" colours
# After :  theme1.xml = "Office Theme" colours, theme2.xml = "Integral" colours
# (fontScheme / fmtScheme are identical between the two themes, only the
#  clrScheme - i.e. the 12 theme colours - and the theme/clrScheme "name"
#  differ.)
#
# The PowerPoint object model only exposes one live ThemeColorScheme (the
# one backing the slide master's Design / theme1.xml), so we repaint that
# scheme with the "Office Theme" palette, which is the half of the swap
# that is reachable through the Design / Theme COM surface.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function Set-ThemeRGB {
    param($Scheme, [int]$Index, [string]$Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    $rgb = $r + ($g * 256) + ($b * 65536)
    $Scheme.Item($Index).RGB = $rgb
}

# Index order exposed by ThemeColorScheme: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
Set-ThemeRGB -Scheme $colorScheme -Index 1  -Hex "000000"   # dk1
Set-ThemeRGB -Scheme $colorScheme -Index 2  -Hex "FFFFFF"   # lt1
Set-ThemeRGB -Scheme $colorScheme -Index 3  -Hex "44546A"   # dk2
Set-ThemeRGB -Scheme $colorScheme -Index 4  -Hex "E7E6E6"   # lt2
Set-ThemeRGB -Scheme $colorScheme -Index 5  -Hex "5B9BD5"   # accent1
Set-ThemeRGB -Scheme $colorScheme -Index 6  -Hex "ED7D31"   # accent2
Set-ThemeRGB -Scheme $colorScheme -Index 7  -Hex "A5A5A5"   # accent3
Set-ThemeRGB -Scheme $colorScheme -Index 8  -Hex "FFC000"   # accent4
Set-ThemeRGB -Scheme $colorScheme -Index 9  -Hex "4472C4"   # accent5
Set-ThemeRGB -Scheme $colorScheme -Index 10 -Hex "70AD47"   # accent6
Set-ThemeRGB -Scheme $colorScheme -Index 11 -Hex "0563C1"   # hlink
Set-ThemeRGB -Scheme $colorScheme -Index 12 -Hex "954F72"   # folHlink
